# This script updates betting-odds values in Sheet1 of the workbook,
# matching the edits described by the provided OOXML diff.
# Only numeric <v> values in columns G..AO for the listed rows change;
# all other cells (labels, styles, row/col layout) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$ws.Range("G3").Value = 1.47  # was 1.42
$ws.Range("H3").Value = 3.65  # was 3.9
$ws.Range("I3").Value = 7.8  # was 8.5
$ws.Range("J3").Value = 2  # was 1.91
$ws.Range("K3").Value = 2.1  # was 2.18
$ws.Range("L3").Value = 7.4  # was 7.8
$ws.Range("M3").Value = 1.42  # was 1.4
$ws.Range("N3").Value = 2.65  # was 2.72
$ws.Range("O3").Value = 2.25  # was 2.18
$ws.Range("P3").Value = 1.57  # was 1.62
$ws.Range("Q3").Value = 3.95  # was 3.8
$ws.Range("R3").Value = 1.21  # was 1.22
$ws.Range("S3").Value = 1.47  # was 1.44
$ws.Range("T3").Value = 2.5  # was 2.57
$ws.Range("U3").Value = 2.4  # was 2.5
$ws.Range("V3").Value = 1.5  # was 1.47
$ws.Range("X3").Value = 5.6  # was 5.3
$ws.Range("Y3").Value = 8.75  # was 9
$ws.Range("Z3").Value = 9.5  # was 8.5
$ws.Range("AC3").Value = 6  # was 6.3
$ws.Range("AD3").Value = 7.5  # was 8
$ws.Range("AE3").Value = 26  # was 28
$ws.Range("AF3").Value = 175  # was 200
$ws.Range("AH3").Value = 14.5  # was 16
$ws.Range("AI3").Value = 50  # was 55
$ws.Range("AJ3").Value = 26  # was 29
$ws.Range("AL3").Value = 120  # was 150
$ws.Range("AM3").Value = 120  # was 150
$ws.Range("AO3").Value = 6  # was 6.3

# Row 6
$ws.Range("J6").Value = 2.67  # was 2.7
$ws.Range("L6").Value = 3.75  # was 3.65
$ws.Range("W6").Value = 6.4  # was 6.2
$ws.Range("X6").Value = 9  # was 8.75
$ws.Range("Y6").Value = 7.3  # was 7.4
$ws.Range("AA6").Value = 14  # was 14.5
$ws.Range("AB6").Value = 21  # was 22
$ws.Range("AC6").Value = 8.25  # was 8.5
$ws.Range("AH6").Value = 7.7  # was 7.9
$ws.Range("AJ6").Value = 9.5  # was 9.25
$ws.Range("AL6").Value = 24  # was 23
$ws.Range("AM6").Value = 29  # was 28

# Row 7
$ws.Range("H7").Value = 3.1  # was 3.15
$ws.Range("S7").Value = 1.42  # was 1.4
$ws.Range("T7").Value = 2.37  # was 2.4
$ws.Range("AC7").Value = 8  # was 8.25

# Row 9
$ws.Range("G9").Value = 1.35  # was 1.36
$ws.Range("H9").Value = 4.55  # was 4.5
$ws.Range("I9").Value = 7.9  # was 7.6
$ws.Range("J9").Value = 1.78  # was 1.82
$ws.Range("K9").Value = 2.45  # was 2.42
$ws.Range("L9").Value = 6.6  # was 6.4
$ws.Range("O9").Value = 1.5  # was 1.52
$ws.Range("P9").Value = 2.25  # was 2.22
$ws.Range("Q9").Value = 2.2  # was 2.22
$ws.Range("R9").Value = 1.53  # was 1.52
$ws.Range("U9").Value = 1.75  # was 1.72
$ws.Range("V9").Value = 1.87  # was 1.88
$ws.Range("W9").Value = 8.25  # was 8
$ws.Range("Z9").Value = 9  # was 9.25
$ws.Range("AD9").Value = 9.5  # was 9.25
$ws.Range("AE9").Value = 17.5  # was 17
$ws.Range("AF9").Value = 70  # was 65
$ws.Range("AJ9").Value = 24  # was 23
$ws.Range("AM9").Value = 65  # was 60

# Row 10
$ws.Range("H10").Value = 3.75  # was 3.8
$ws.Range("I10").Value = 5.5  # was 5.4
$ws.Range("J10").Value = 2.15  # was 2.12
$ws.Range("K10").Value = 2.15  # was 2.18
$ws.Range("L10").Value = 5.4  # was 5.3
$ws.Range("W10").Value = 6.7  # was 6.8
$ws.Range("AB10").Value = 25  # was 26
$ws.Range("AC10").Value = 10.75  # was 11
$ws.Range("AD10").Value = 7.4  # was 7.5
$ws.Range("AF10").Value = 70  # was 75
$ws.Range("AG10").Value = 500  # was 600
$ws.Range("AH10").Value = 16  # was 15.5

# Row 11
$ws.Range("U11").Value = 1.76  # was 1.78
$ws.Range("V11").Value = 1.95  # was 1.93

# Row 12
$ws.Range("G12").Value = 2.02  # was 1.93
$ws.Range("I12").Value = 3.1  # was 3.35
$ws.Range("J12").Value = 2.6  # was 2.5
$ws.Range("L12").Value = 3.55  # was 3.8
$ws.Range("M12").Value = 1.23  # was 1.24
$ws.Range("N12").Value = 3.35  # was 3.3
$ws.Range("O12").Value = 1.7  # was 1.72
$ws.Range("P12").Value = 1.91  # was 1.88
$ws.Range("Q12").Value = 2.65  # was 2.7
$ws.Range("R12").Value = 1.36  # was 1.35
$ws.Range("U12").Value = 1.62  # was 1.65
$ws.Range("V12").Value = 2.02  # was 1.98
$ws.Range("W12").Value = 8.5  # was 8
$ws.Range("X12").Value = 10.5  # was 9.5
$ws.Range("Y12").Value = 8.75  # was 8.5
$ws.Range("Z12").Value = 18.5  # was 16.5
$ws.Range("AA12").Value = 15.5  # was 15
$ws.Range("AB12").Value = 24  # was 25
$ws.Range("AC12").Value = 12  # was 11.5
$ws.Range("AE12").Value = 13.5  # was 14
$ws.Range("AF12").Value = 55  # was 60
$ws.Range("AG12").Value = 400  # was 450
$ws.Range("AH12").Value = 10.75  # was 11.25
$ws.Range("AI12").Value = 17  # was 18.5
$ws.Range("AJ12").Value = 11  # was 11.75
$ws.Range("AK12").Value = 40  # was 45
$ws.Range("AL12").Value = 25  # was 28
$ws.Range("AM12").Value = 32  # was 35

# Row 13
$ws.Range("O13").Value = 2.08  # was 2.05
$ws.Range("P13").Value = 1.73  # was 1.75
$ws.Range("U13").Value = 2.1  # was 2
$ws.Range("V13").Value = 1.63  # was 1.73
$ws.Range("AB13").Value = 34  # was 29
$ws.Range("AC13").Value = 9  # was 9.5
$ws.Range("AN13").Value = 1.03  # was 1.04
$ws.Range("AO13").Value = 9  # was 8

# Row 14
$ws.Range("G14").Value = 3.65  # was 3.6
$ws.Range("L14").Value = 2.72  # was 2.75
$ws.Range("M14").Value = 1.44  # was 1.45
$ws.Range("O14").Value = 2.3  # was 2.32
$ws.Range("Q14").Value = 4.05  # was 4.1
$ws.Range("S14").Value = 1.47  # was 1.5
$ws.Range("T14").Value = 2.5  # was 2.42
$ws.Range("U14").Value = 1.98  # was 2
$ws.Range("V14").Value = 1.75  # was 1.72
$ws.Range("AA14").Value = 40  # was 37
$ws.Range("AE14").Value = 16.5  # was 17
$ws.Range("AH14").Value = 6  # was 5.9
$ws.Range("AL14").Value = 19.5  # was 20
$ws.Range("AM14").Value = 35  # was 37

# Row 15
$ws.Range("G15").Value = 10  # was 7.3
$ws.Range("H15").Value = 4.7  # was 4.4
$ws.Range("I15").Value = 1.3  # was 1.4
$ws.Range("J15").Value = 7.9  # was 6.5
$ws.Range("K15").Value = 2.4  # was 2.32
$ws.Range("L15").Value = 1.78  # was 1.91
$ws.Range("M15").Value = 1.21  # was 1.23
$ws.Range("N15").Value = 3.95  # was 3.75
$ws.Range("O15").Value = 1.65  # was 1.7
$ws.Range("P15").Value = 2.12  # was 2.02
$ws.Range("Q15").Value = 2.52  # was 2.67
$ws.Range("R15").Value = 1.47  # was 1.42
$ws.Range("S15").Value = 1.34  # was 1.36
$ws.Range("T15").Value = 3  # was 2.92
$ws.Range("U15").Value = 2.02  # was 1.93
$ws.Range("V15").Value = 1.72  # was 1.78
$ws.Range("W15").Value = 26  # was 19.5
$ws.Range("X15").Value = 80  # was 50
$ws.Range("Y15").Value = 30  # was 22
$ws.Range("Z15").Value = 300  # was 175
$ws.Range("AA15").Value = 120  # was 80
$ws.Range("AB15").Value = 90  # was 70
$ws.Range("AC15").Value = 8.5  # was 8.25
$ws.Range("AD15").Value = 9.5  # was 8.75
$ws.Range("AE15").Value = 21  # was 19.5
$ws.Range("AF15").Value = 100  # was 90
$ws.Range("AG15").Value = 800  # was 700
$ws.Range("AH15").Value = 6.9  # was 6.8
$ws.Range("AI15").Value = 6.2  # was 6.5
$ws.Range("AJ15").Value = 8.5  # was 8.25
$ws.Range("AK15").Value = 8  # was 9
$ws.Range("AL15").Value = 11  # was 11.5
$ws.Range("AM15").Value = 28  # was 27
$ws.Range("AO15").Value = 8.5  # was 8.25

# Row 16
$ws.Range("M16").Value = 1.67  # was 1.63
$ws.Range("R16").Value = 1.11  # was 1.08
$ws.Range("U16").Value = 2.38  # was 2.25
$ws.Range("V16").Value = 1.53  # was 1.57
$ws.Range("AN16").Value = 1.14  # was 1.11

# Row 17
$ws.Range("O17").Value = 1.95  # was 1.9
$ws.Range("P17").Value = 1.85  # was 1.9
$ws.Range("Q17").Value = 3.4  # was 3.25
$ws.Range("R17").Value = 1.3  # was 1.33

# Row 18
$ws.Range("O18").Value = 2.05  # was 2.08
$ws.Range("P18").Value = 1.75  # was 1.73
$ws.Range("Q18").Value = 3.5  # was 3.75
$ws.Range("R18").Value = 1.29  # was 1.25
$ws.Range("AN18").Value = 1.06  # was 1.07
$ws.Range("AO18").Value = 10  # was 9

# Row 19
$ws.Range("AN19").Value = 1.06  # was 1.07
$ws.Range("AO19").Value = 10  # was 9

